$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("583:583").Insert()

$ws.Range("A583").Value = 10
$ws.Range("B583").Value = "Vega Modelo de Temuco"
$ws.Range("C583").Value = "La Araucanía"
$ws.Range("D583").Value = 45212
$ws.Range("E583").Value = 9
$ws.Range("F583").Value = 100112024
$ws.Range("G583").Value = "Choclo"
$ws.Range("H583").Value = "Choclero"
$ws.Range("I583").Value = "Primera"
$ws.Range("J583").Value = 20
$ws.Range("K583").Value = 22000
$ws.Range("L583").Value = 22000
$ws.Range("M583").Value = 22000
$ws.Range("N583").Value = "$/malla 50 unidades"
$ws.Range("O583").Value = "Región de Arica y Parinacota"
$ws.Range("P583").Value = 440
$ws.Range("Q583").Value = 50
$ws.Range("R583").Value = "Hortaliza"
